# Key Personnel List - update the placeholder/instruction text above the
# personnel table. The "Full Title" line becomes "LONG" and the "Address"
# line becomes "ADDRESS" (these act as markers used by the PDF-merge
# automation referenced in the commit message).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value() = "LONG"
$ws.Range("A3").Value() = "ADDRESS"
